# Applies the Coinranking crypto price/volume refresh for Tue Jan 30 19:33:02 UTC 2024 run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.432.63"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "2.372.36"
$ws.Range("E3").Value = "  +3.13%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'309.45"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").Value = "'105.06"
$ws.Range("E6").Value = "  +4.51%  "
$ws.Range("D7").Value = "'0.518"
$ws.Range("E7").Value = "  -3.34%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("D10").Value = "'36.20"
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("D11").Value = "'53.34"
$ws.Range("E11").Value = "  +2.40%  "
$ws.Range("E12").Value = "  -1.54%  "
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("E14").Value = "  -2.27%  "
$ws.Range("D15").Value = "2.737.75"
$ws.Range("E15").Value = "  +3.05%  "
$ws.Range("D16").Value = "'15.60"
$ws.Range("E16").Value = "  +4.11%  "
$ws.Range("D17").Value = "2.372.49"
$ws.Range("E17").Value = "  +2.83%  "
$ws.Range("D18").Value = "'0.813"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").Value = "43.386.52"
$ws.Range("E19").Value = "  +0.89%  "
$ws.Range("D20").Value = "'12.04"
$ws.Range("E20").Value = "  -3.53%  "
$ws.Range("D21").Value = "'6.31"
$ws.Range("E21").Value = "  +3.88%  "
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("D23").Value = "'68.26"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").Value = "'241.57"
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("E25").Value = "  +1.66%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("E28").Value = "  +5.93%  "
$ws.Range("D29").Value = "'3.85"
$ws.Range("E29").Value = "  -3.76%  "
$ws.Range("D30").Value = "'2.30"
$ws.Range("E30").Value = "  +8.63%  "
$ws.Range("D31").Value = "'36.78"
$ws.Range("E31").Value = "  -4.26%  "
$ws.Range("D32").Value = "'9.56"
$ws.Range("E32").Value = "  -0.91%  "
$ws.Range("D33").Value = "'161.85"
$ws.Range("E33").Value = "  -3.14%  "
$ws.Range("E34").Value = "  -0.96%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").Value = "'18.39"
$ws.Range("E36").Value = "  +3.57%  "
$ws.Range("E37").Value = "  +6.16%  "
$ws.Range("E38").Value = "  +12.95%  "
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("D40").Value = "'0.0743"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("E41").Value = "  +5.75%  "
$ws.Range("E42").Value = "  +0.22%  "
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("D44").Value = "'2.67"
$ws.Range("E44").Value = "  +17.39%  "
$ws.Range("D45").Value = "'19.91"
$ws.Range("E45").Value = "  +3.96%  "
$ws.Range("D46").Value = "2.005.78"
$ws.Range("E46").Value = "  +1.86%  "
$ws.Range("E47").Value = "  +0.66%  "
$ws.Range("D48").Value = "'3.14"
$ws.Range("E48").Value = "  +4.01%  "
$ws.Range("D49").Value = "'10.63"
$ws.Range("E49").Value = "  +8.02%  "
$ws.Range("D50").Value = "'58.17"
$ws.Range("E50").Value = "  +4.73%  "
$ws.Range("E51").Value = "  +1.11%  "
